# Adds scraping support for extra bowling attributes:
#  - clears out leftover empty placeholder cells on "ODI Batting Extra"
#  - adds a new "ODI Bowling Extra" worksheet with MATCH_CODE / MAIDEN_OVERS /
#    PERCENT_WICKETS_OF_ALL columns

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Clean up empty placeholder cells left in "ODI Batting Extra"
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")

$emptyCellRefs = @(
    "C2","D2","E2",
    "B7","C7","D7","E7",
    "C13","D13","E13",
    "E14",
    "B15","C15","D15","E15",
    "E16",
    "B17","C17","D17","E17",
    "B19","C19","D19","E19",
    "B20","C20","D20","E20",
    "B21","C21","D21","E21"
)

foreach ($ref in $emptyCellRefs) {
    $battingExtra.Range($ref).ClearContents()
}

# ---------------------------------------------------------------------------
# 2) Add the new "ODI Bowling Extra" worksheet after "ODI Batting Extra"
# ---------------------------------------------------------------------------
$bowlingExtra = $wb.Worksheets.Add($null, $battingExtra)
$bowlingExtra.Name = "ODI Bowling Extra"

# Header row
$bowlingExtra.Range("A1").Value = "MATCH_CODE"
$bowlingExtra.Range("B1").Value = "MAIDEN_OVERS"
$bowlingExtra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# Match the bold/bordered/centered header style used on the other sheets
$battingExtra.Range("A1").Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122)

# Data rows - force text formatting so numeric-looking strings (e.g. "0",
# "10.00%") are preserved exactly as text rather than being reinterpreted.
$bowlingExtra.Range("A2:C21").NumberFormat = "@"

$matchCodes   = @("3701","3702","3754","3765","3774","3780","3784","3799","3813","3857","3938","3971","4031","4047","4083","4086","4250","4251","4296","4357")
$maidenOvers  = @("0","0","0","0","","0","0","0","0","0","","0","0","0","0","0","0","0","0","")
$pctWicketsAll = @("","","","","","","","","","10.00%","","","10.00%","","","","","","","")

for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $row = $i + 2
    $bowlingExtra.Cells.Item($row, 1).Value = $matchCodes[$i]
    $bowlingExtra.Cells.Item($row, 2).Value = $maidenOvers[$i]
    $bowlingExtra.Cells.Item($row, 3).Value = $pctWicketsAll[$i]
}
